# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# for various leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# reflecting refreshed market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 895342.1
$ws.Range("I15").Value = 895342.1
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2686026.3
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2685857.3
$ws.Range("H70").Value = 2730
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 3825
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 11475
$ws.Range("M70").Value = -5730
$ws.Range("H73").Value = 2730
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 3825
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 11475
$ws.Range("M73").Value = -5064
$ws.Range("H100").Value = 2224.0908
$ws.Range("I100").Value = 1507.8572
$ws.Range("J100").Value = 3477.5
$ws.Range("K100").Value = 1507.8572
$ws.Range("L100").Value = 3477.5
$ws.Range("M100").Value = -966.8571999999999
$ws.Range("N100").Value = -4559.5
$ws.Range("H106").Value = 29335330
$ws.Range("I106").Value = 36668184
$ws.Range("J106").Value = 3902
$ws.Range("K106").Value = 36668184
$ws.Range("L106").Value = 3902
$ws.Range("M106").Value = -36667553
$ws.Range("H132").Value = 3681.139
$ws.Range("I132").Value = 3430.0645
$ws.Range("J132").Value = 5237.8
$ws.Range("K132").Value = 10290.1935
$ws.Range("L132").Value = 15713.4
$ws.Range("M132").Value = -7760.193499999999
$ws.Range("N132").Value = -20773.4
$ws.Range("H137").Value = 7526.3076
$ws.Range("I137").Value = 6129.722
$ws.Range("J137").Value = 8723.380999999999
$ws.Range("K137").Value = 18389.166
$ws.Range("L137").Value = 26170.143
$ws.Range("M137").Value = -15839.166
$ws.Range("N137").Value = -31270.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4705.1772
$ws.Range("I32").Value = 3097.0679
$ws.Range("J32").Value = 36331.332
$ws.Range("K32").Value = 3097.0679
$ws.Range("L32").Value = 36331.332
$ws.Range("M32").Value = -2810.0679
$ws.Range("N32").Value = -36905.332
$ws.Range("H74").Value = 4443.931
$ws.Range("I74").Value = 2415.9473
$ws.Range("J74").Value = 8297.1
$ws.Range("K74").Value = 2415.9473
$ws.Range("L74").Value = 8297.1
$ws.Range("M74").Value = -1541.9473
$ws.Range("H77").Value = 4443.931
$ws.Range("I77").Value = 2415.9473
$ws.Range("J77").Value = 8297.1
$ws.Range("K77").Value = 12079.7365
$ws.Range("L77").Value = 41485.5
$ws.Range("M77").Value = -7711.736499999999
$ws.Range("H102").Value = 1353.1428
$ws.Range("I102").Value = 498.25
$ws.Range("J102").Value = 2493
$ws.Range("K102").Value = 498.25
$ws.Range("L102").Value = 2493
$ws.Range("M102").Value = 1123.75
$ws.Range("H110").Value = 71431520
$ws.Range("I110").Value = 3442.5
$ws.Range("J110").Value = 500000000
$ws.Range("K110").Value = 3442.5
$ws.Range("L110").Value = 500000000
$ws.Range("M110").Value = -1397.5
$ws.Range("N110").Value = -500004090
$ws.Range("H122").Value = 7466.3335
$ws.Range("I122").Value = 7466.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 22399.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -19949.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2322.889
$ws.Range("I99").Value = 1499
$ws.Range("J99").Value = 3352.75
$ws.Range("K99").Value = 1499
$ws.Range("L99").Value = 3352.75
$ws.Range("M99").Value = -1
$ws.Range("H107").Value = 35714790
$ws.Range("I107").Value = 388.33334
$ws.Range("J107").Value = 100000720
$ws.Range("K107").Value = 388.33334
$ws.Range("L107").Value = 100000720
$ws.Range("M107").Value = 1531.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1422
$ws.Range("I16").Value = 1152.5
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1152.5
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -865.5
$ws.Range("N16").Value = -3074
$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20470
$ws.Range("M21").ClearContents()
$ws.Range("H31").Value = 4871.4346
$ws.Range("I31").Value = 3178.3022
$ws.Range("J31").Value = 7671.615
$ws.Range("K31").Value = 3178.3022
$ws.Range("L31").Value = 7671.615
$ws.Range("M31").Value = -2883.3022
$ws.Range("N31").Value = -8261.615
$ws.Range("H34").Value = 4871.4346
$ws.Range("I34").Value = 3178.3022
$ws.Range("J34").Value = 7671.615
$ws.Range("K34").Value = 3178.3022
$ws.Range("L34").Value = 7671.615
$ws.Range("M34").Value = -2976.3022
$ws.Range("N34").Value = -8075.615
$ws.Range("H58").Value = 9465.267
$ws.Range("I58").Value = 4985.5713
$ws.Range("J58").Value = 13385
$ws.Range("K58").Value = 4985.5713
$ws.Range("L58").Value = 13385
$ws.Range("M58").Value = -4782.5713
$ws.Range("N58").Value = -13791
$ws.Range("H60").Value = 25558.4
$ws.Range("I60").Value = 4500
$ws.Range("J60").Value = 39597.332
$ws.Range("K60").Value = 4500
$ws.Range("L60").Value = 39597.332
$ws.Range("M60").Value = -3989
$ws.Range("N60").Value = -40619.332
$ws.Range("H105").Value = 52633028
$ws.Range("I105").Value = 1529.6111
$ws.Range("J105").Value = 1000000000
$ws.Range("K105").Value = 1529.6111
$ws.Range("L105").Value = 1000000000
$ws.Range("M105").Value = 217.3888999999999
$ws.Range("H113").Value = 1422
$ws.Range("I113").Value = 1152.5
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1152.5
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 1017.5
$ws.Range("N113").Value = -6840
$ws.Range("H132").Value = 4949.2354
$ws.Range("I132").Value = 3801.6428
$ws.Range("J132").Value = 10304.667
$ws.Range("K132").Value = 11404.9284
$ws.Range("L132").Value = 30914.001
$ws.Range("M132").Value = -8874.928400000001
$ws.Range("N132").Value = -35974.001
$ws.Range("H134").Value = 7313.5
$ws.Range("I134").Value = 7313.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 21940.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -19405.5
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 9465.267
$ws.Range("I136").Value = 4985.5713
$ws.Range("J136").Value = 13385
$ws.Range("K136").Value = 14956.7139
$ws.Range("L136").Value = 40155
$ws.Range("M136").Value = -12406.7139
$ws.Range("N136").Value = -45255

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3403.889
$ws.Range("I68").Value = 2335.7778
$ws.Range("J68").Value = 3937.9443
$ws.Range("K68").Value = 7007.3334
$ws.Range("L68").Value = 11813.8329
$ws.Range("M68").Value = -6196.3334
$ws.Range("N68").Value = -13435.8329
$ws.Range("H70").Value = 15084.875
$ws.Range("I70").Value = 10169.75
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 30509.25
$ws.Range("L70").Value = 60000
$ws.Range("M70").Value = -30194.25
$ws.Range("H71").Value = 3403.889
$ws.Range("I71").Value = 2335.7778
$ws.Range("J71").Value = 3937.9443
$ws.Range("K71").Value = 21022.0002
$ws.Range("L71").Value = 35441.4987
$ws.Range("M71").Value = -16966.0002
$ws.Range("N71").Value = -43553.4987
$ws.Range("H73").Value = 15084.875
$ws.Range("I73").Value = 10169.75
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 30509.25
$ws.Range("L73").Value = 60000
$ws.Range("M73").Value = -29417.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3418.35
$ws.Range("I102").Value = 1722.6364
$ws.Range("J102").Value = 5490.8887
$ws.Range("K102").Value = 1722.6364
$ws.Range("L102").Value = 5490.8887
$ws.Range("M102").Value = -100.6364000000001
$ws.Range("H113").Value = 2573.5
$ws.Range("I113").Value = 1851.9546
$ws.Range("J113").Value = 4160.9
$ws.Range("K113").Value = 1851.9546
$ws.Range("L113").Value = 4160.9
$ws.Range("M113").Value = 318.0454
$ws.Range("N113").Value = -8500.9
$ws.Range("H122").Value = 4997.8335
$ws.Range("I122").Value = 4997.8335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14993.5005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12543.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13792
$ws.Range("I7").Value = 13792
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 13792
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -13680
$ws.Range("N7").ClearContents()
$ws.Range("H68").Value = 5586.6313
$ws.Range("I68").Value = 5509.125
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 5509.125
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -4760.125
$ws.Range("N68").Value = -7498
$ws.Range("H71").Value = 5586.6313
$ws.Range("I71").Value = 5509.125
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 27545.625
$ws.Range("L71").Value = 30000
$ws.Range("M71").Value = -23801.625
$ws.Range("N71").Value = -37488
$ws.Range("H93").Value = 1693.1818
$ws.Range("I93").Value = 1693.1818
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1693.1818
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -445.1818000000001
$ws.Range("H97").Value = 129292.375
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 129292.375
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 129292.375
$ws.Range("N97").Value = -131274.375
$ws.Range("H126").Value = 13792
$ws.Range("I126").Value = 13792
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 41376
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -38906
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 516.3333
$ws.Range("I100").Value = 374.6
$ws.Range("J100").Value = 1225
$ws.Range("K100").Value = 749.2
$ws.Range("L100").Value = 2450
$ws.Range("M100").Value = -208.2
$ws.Range("N100").Value = -3532
$ws.Range("H132").Value = 3626.1924
$ws.Range("I132").Value = 2892.0454
$ws.Range("J132").Value = 7664
$ws.Range("K132").Value = 8676.136200000001
$ws.Range("L132").Value = 22992
$ws.Range("M132").Value = -6146.136200000001
$ws.Range("N132").Value = -28052
$ws.Range("H135").Value = 165000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 165000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 165000
$ws.Range("N135").Value = -175140
$ws.Range("H136").Value = 6516.604
$ws.Range("I136").Value = 5961.5938
$ws.Range("J136").Value = 7362.3335
$ws.Range("K136").Value = 17884.7814
$ws.Range("L136").Value = 22087.0005
$ws.Range("M136").Value = -15334.7814
$ws.Range("N136").Value = -27187.0005

